$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.588.28'
$ws.Range("E2").Value = '  +1.90%  '
$ws.Range("D3").Value = '2.550.43'
$ws.Range("E3").Value = '  +5.04%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.47'
$ws.Range("E5").Value = '  +2.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.08'
$ws.Range("E6").Value = '  +8.99%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").Value = '2.548.64'
$ws.Range("E9").Value = '  +5.05%  '
$ws.Range("E10").Value = '  +2.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.75'
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.359'
$ws.Range("E13").Value = '  +3.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.53'
$ws.Range("E14").Value = '  +9.04%  '
$ws.Range("D15").Value = '3.006.81'
$ws.Range("E15").Value = '  +5.11%  '
$ws.Range("D16").Value = '63.495.59'
$ws.Range("E16").Value = '  +1.91%  '
$ws.Range("D18").Value = '2.540.80'
$ws.Range("E18").Value = '  +4.70%  '
$ws.Range("E19").Value = '  +4.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.18'
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("E21").Value = '  +4.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.86'
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.14'
$ws.Range("E24").Value = '  +1.55%  '
$ws.Range("E25").Value = '  -0.68%  '
$ws.Range("E26").Value = '  +5.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.54'
$ws.Range("E27").Value = '  +14.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.60'
$ws.Range("E28").Value = '  +5.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("E30").Value = '  +12.84%  '
$ws.Range("D31").Value = '0.0₃0833'
$ws.Range("E31").Value = '  +6.33%  '
$ws.Range("E32").Value = '  +4.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '177.99'
$ws.Range("E33").Value = '  +3.37%  '
$ws.Range("E34").Value = '  +9.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '424.27'
$ws.Range("E35").Value = '  +12.33%  '
$ws.Range("E36").Value = '  +2.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.21'
$ws.Range("E37").Value = '  +3.59%  '
$ws.Range("E38").Value = '  +0.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.80'
$ws.Range("E39").Value = '  +6.81%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.80'
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '154.11'
$ws.Range("E43").Value = '  +6.37%  '
$ws.Range("E44").Value = '  +4.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.13'
$ws.Range("E45").Value = '  +2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.612'
$ws.Range("E46").Value = '  +3.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0532'
$ws.Range("E47").Value = '  +2.63%  '
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("E49").Value = '  +7.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.71'
$ws.Range("E50").Value = '  +4.58%  '
$ws.Range("E51").Value = '  +8.06%  '
